$p = $ppt.ActivePresentation

# 1) Remove the trailing slides (9-24) that are no longer part of the deck.
for ($i = $p.Slides.Count; $i -ge 9; $i--) {
    $p.Slides.Item($i).Delete()
}

# 2) Slide 1 - title slide subtitle becomes the new LSTM paper intro block.
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "Word-Level LSTM Model for Sentence Completion`rusing Shakespeare’s Plays`rAfaq Alam`rB.Sc Data Science`rNUCES-FAST`rIslamabad, Pakistan`rAbstract—This paper presents a word-level Long Short-Term`rMemory (LSTM) model trained on Shakespeare’s plays to predict`rthe next word in a sequence. The model is trained using`rTensorFlow and Keras on a dataset containing Shakespearean`rdialogues"

# 3) Slide 2 - Key Points bullets.
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(2).TextFrame.TextRange.Text = "A word-level LSTM model is trained on Shakespeare’s plays to predict the next word in a sequence`rThe model is integrated with a user-friendly interface that provides real-time word suggestions"

# 4) Slide 3 - single bullet.
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(2).TextFrame.TextRange.Text = "The report discusses the preprocessing steps, model architecture, results, and challenges encountered during implementation ."

# 5) Slide 4 - two bullets.
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(2).TextFrame.TextRange.Text = "The Shakespeare Plays dataset was obtained from Kaggle and includes dialogues from various plays written by William Shakespeare`rThe study also explores how hyperparameters are used to evaluate the accuracy of the model’s accuracy and evaluates the coherence of the sentences ."

# 6) Slide 5 - title goes from blank to "Key Points"; bullets change.
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Key Points"
$s5.Shapes.Item(2).TextFrame.TextRange.Text = "The e300 b128 model achieved 91.14% accuracy`rLarger batch sizes (e.g., 128) lead to faster convergence but may generalize less than smaller batch sizes"

# 7) Slide 6 - title goes from "Key Points" to blank; bullets collapse to one.
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = ""
$s6.Shapes.Item(2).TextFrame.TextRange.Text = "Increasing the number of epochs improves accuracy, as seen in the e300 model ."

# 8) Slide 7 - title goes from blank to "Key Points"; bullets grow to two.
$s7 = $p.Slides.Item(7)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Key Points"
$s7.Shapes.Item(2).TextFrame.TextRange.Text = "The e250 b64 model achieves a balance between accu-phthalracy and loss, making it a viable alternative to the highest-performing model`rA higher number of training epochs and a moderately-high batch size improve model performance, but diminishing returns can occur after a certain threshold ."

# 9) Slide 8 - title goes from "Key Points" to "Summary"; bullets collapse to one merged summary.
$s8 = $p.Slides.Item(8)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Summary"
$s8.Shapes.Item(2).TextFrame.TextRange.Text = " The e250 b64 model achieves a balance between accu-phthalracy and loss, making it a viable alternative to the highest-performing model . A higher number of training epochs and a moderately-high batch size improve model performance, but diminishing returns can occur after a certain threshold ."

Write-Output ("Final slide count: " + $p.Slides.Count)
